{"js": "/*\n * Apply the table/date-line text updates described by the commit diff.\n * Each table cell's math expression and the date paragraph are replaced\n * with their new value, in document order (row-major for the table).\n * All run/paragraph formatting is left untouched.\n */\n\n// New text for the date paragraph (first paragraph in the body).\nconst NEW_DATE = \"2024-11-29 Friday\";\n\n// New text for every cell of the (single) table, row-major, 20 rows x 5 cols.\nconst NEW_TABLE_VALUES = [\n  [\n    \"94-26=\",\n    \"97-28=\",\n    \"30-1=\",\n    \"44+36=\",\n    \"42+20=\"\n  ],\n  [\n    \"89-33=\",\n    \"43+49=\",\n    \"28-15=\",\n    \"42-39=\",\n    \"58+30=\"\n  ],\n  [\n    \"86-50=\",\n    \"84+6=\",\n    \"21+12=\",\n    \"43+14=\",\n    \"20+29=\"\n  ],\n  [\n    \"4+3=\",\n    \"84-0=\",\n    \"34-25=\",\n    \"85-42=\",\n    \"78-67=\"\n  ],\n  [\n    \"52+28=\",\n    \"90-83=\",\n    \"89-68=\",\n    \"30+7=\",\n    \"37-27=\"\n  ],\n  [\n    \"20+35=\",\n    \"67-10=\",\n    \"47+44=\",\n    \"96-66=\",\n    \"54+13=\"\n  ],\n  [\n    \"78-5=\",\n    \"49+12=\",\n    \"80-4=\",\n    \"15+36=\",\n    \"31+44=\"\n  ],\n  [\n    \"90-61=\",\n    \"79+17=\",\n    \"81-29=\",\n    \"80+3=\",\n    \"79-39=\"\n  ],\n  [\n    \"32-18=\",\n    \"7+85=\",\n    \"51+25=\",\n    \"12-5=\",\n    \"22+32=\"\n  ],\n  [\n    \"32+22=\",\n    \"31+43=\",\n    \"42+11=\",\n    \"34+25=\",\n    \"46-13=\"\n  ],\n  [\n    \"16+80=\",\n    \"45+35=\",\n    \"15-15=\",\n    \"49-17=\",\n    \"90-21=\"\n  ],\n  [\n    \"32-5=\",\n    \"87-70=\",\n    \"40+46=\",\n    \"43-29=\",\n    \"28+71=\"\n  ],\n  [\n    \"1+95=\",\n    \"13+25=\",\n    \"86-42=\",\n    \"36-19=\",\n    \"85-51=\"\n  ],\n  [\n    \"83+0=\",\n    \"18+44=\",\n    \"75-47=\",\n    \"5+32=\",\n    \"68-52=\"\n  ],\n  [\n    \"56+29=\",\n    \"18+31=\",\n    \"15+68=\",\n    \"12-0=\",\n    \"55+35=\"\n  ],\n  [\n    \"55+24=\",\n    \"58+1=\",\n    \"53-40=\",\n    \"85-7=\",\n    \"77-66=\"\n  ],\n  [\n    \"46-46=\",\n    \"4+25=\",\n    \"36+42=\",\n    \"59+28=\",\n    \"56+24=\"\n  ],\n  [\n    \"74-54=\",\n    \"12+41=\",\n    \"10+19=\",\n    \"85-13=\",\n    \"28-3=\"\n  ],\n  [\n    \"14+53=\",\n    \"80-49=\",\n    \"53+0=\",\n    \"54-29=\",\n    \"35-21=\"\n  ],\n  [\n    \"32+64=\",\n    \"18+7=\",\n    \"0+31=\",\n    \"61-25=\",\n    \"21+48=\"\n  ]\n];\n\nconst body = context.document.body;\n\n// --- 1. Update the date/weekday line -------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  paragraphs.items[0].insertText(NEW_DATE, Word.InsertLocation.replace);\n}\n\n// --- 2. Update every arithmetic-problem cell in the table -----------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.load(\"rowCount, values\");\n  await context.sync();\n\n  // Assigning a full 2-D array to `.values` rewrites each cell's text\n  // while preserving the existing run/paragraph formatting.\n  table.values = NEW_TABLE_VALUES;\n}\n\nawait context.sync();\n", "ps1": "# Update the date/weekday line (first paragraph in the document body).\n$d = $word.ActiveDocument\n$d.Paragraphs.Item(1).Range.Text = \"2024-11-29 Friday\"\n\n# Update every arithmetic-problem cell in the (single) table, row-major.\n$t = $d.Tables.Item(1)\n$t.Cell(1, 1).Range.Text = \"94-26=\"\n$t.Cell(1, 2).Range.Text = \"97-28=\"\n$t.Cell(1, 3).Range.Text = \"30-1=\"\n$t.Cell(1, 4).Range.Text = \"44+36=\"\n$t.Cell(1, 5).Range.Text = \"42+20=\"\n$t.Cell(2, 1).Range.Text = \"89-33=\"\n$t.Cell(2, 2).Range.Text = \"43+49=\"\n$t.Cell(2, 3).Range.Text = \"28-15=\"\n$t.Cell(2, 4).Range.Text = \"42-39=\"\n$t.Cell(2, 5).Range.Text = \"58+30=\"\n$t.Cell(3, 1).Range.Text = \"86-50=\"\n$t.Cell(3, 2).Range.Text = \"84+6=\"\n$t.Cell(3, 3).Range.Text = \"21+12=\"\n$t.Cell(3, 4).Range.Text = \"43+14=\"\n$t.Cell(3, 5).Range.Text = \"20+29=\"\n$t.Cell(4, 1).Range.Text = \"4+3=\"\n$t.Cell(4, 2).Range.Text = \"84-0=\"\n$t.Cell(4, 3).Range.Text = \"34-25=\"\n$t.Cell(4, 4).Range.Text = \"85-42=\"\n$t.Cell(4, 5).Range.Text = \"78-67=\"\n$t.Cell(5, 1).Range.Text = \"52+28=\"\n$t.Cell(5, 2).Range.Text = \"90-83=\"\n$t.Cell(5, 3).Range.Text = \"89-68=\"\n$t.Cell(5, 4).Range.Text = \"30+7=\"\n$t.Cell(5, 5).Range.Text = \"37-27=\"\n$t.Cell(6, 1).Range.Text = \"20+35=\"\n$t.Cell(6, 2).Range.Text = \"67-10=\"\n$t.Cell(6, 3).Range.Text = \"47+44=\"\n$t.Cell(6, 4).Range.Text = \"96-66=\"\n$t.Cell(6, 5).Range.Text = \"54+13=\"\n$t.Cell(7, 1).Range.Text = \"78-5=\"\n$t.Cell(7, 2).Range.Text = \"49+12=\"\n$t.Cell(7, 3).Range.Text = \"80-4=\"\n$t.Cell(7, 4).Range.Text = \"15+36=\"\n$t.Cell(7, 5).Range.Text = \"31+44=\"\n$t.Cell(8, 1).Range.Text = \"90-61=\"\n$t.Cell(8, 2).Range.Text = \"79+17=\"\n$t.Cell(8, 3).Range.Text = \"81-29=\"\n$t.Cell(8, 4).Range.Text = \"80+3=\"\n$t.Cell(8, 5).Range.Text = \"79-39=\"\n$t.Cell(9, 1).Range.Text = \"32-18=\"\n$t.Cell(9, 2).Range.Text = \"7+85=\"\n$t.Cell(9, 3).Range.Text = \"51+25=\"\n$t.Cell(9, 4).Range.Text = \"12-5=\"\n$t.Cell(9, 5).Range.Text = \"22+32=\"\n$t.Cell(10, 1).Range.Text = \"32+22=\"\n$t.Cell(10, 2).Range.Text = \"31+43=\"\n$t.Cell(10, 3).Range.Text = \"42+11=\"\n$t.Cell(10, 4).Range.Text = \"34+25=\"\n$t.Cell(10, 5).Range.Text = \"46-13=\"\n$t.Cell(11, 1).Range.Text = \"16+80=\"\n$t.Cell(11, 2).Range.Text = \"45+35=\"\n$t.Cell(11, 3).Range.Text = \"15-15=\"\n$t.Cell(11, 4).Range.Text = \"49-17=\"\n$t.Cell(11, 5).Range.Text = \"90-21=\"\n$t.Cell(12, 1).Range.Text = \"32-5=\"\n$t.Cell(12, 2).Range.Text = \"87-70=\"\n$t.Cell(12, 3).Range.Text = \"40+46=\"\n$t.Cell(12, 4).Range.Text = \"43-29=\"\n$t.Cell(12, 5).Range.Text = \"28+71=\"\n$t.Cell(13, 1).Range.Text = \"1+95=\"\n$t.Cell(13, 2).Range.Text = \"13+25=\"\n$t.Cell(13, 3).Range.Text = \"86-42=\"\n$t.Cell(13, 4).Range.Text = \"36-19=\"\n$t.Cell(13, 5).Range.Text = \"85-51=\"\n$t.Cell(14, 1).Range.Text = \"83+0=\"\n$t.Cell(14, 2).Range.Text = \"18+44=\"\n$t.Cell(14, 3).Range.Text = \"75-47=\"\n$t.Cell(14, 4).Range.Text = \"5+32=\"\n$t.Cell(14, 5).Range.Text = \"68-52=\"\n$t.Cell(15, 1).Range.Text = \"56+29=\"\n$t.Cell(15, 2).Range.Text = \"18+31=\"\n$t.Cell(15, 3).Range.Text = \"15+68=\"\n$t.Cell(15, 4).Range.Text = \"12-0=\"\n$t.Cell(15, 5).Range.Text = \"55+35=\"\n$t.Cell(16, 1).Range.Text = \"55+24=\"\n$t.Cell(16, 2).Range.Text = \"58+1=\"\n$t.Cell(16, 3).Range.Text = \"53-40=\"\n$t.Cell(16, 4).Range.Text = \"85-7=\"\n$t.Cell(16, 5).Range.Text = \"77-66=\"\n$t.Cell(17, 1).Range.Text = \"46-46=\"\n$t.Cell(17, 2).Range.Text = \"4+25=\"\n$t.Cell(17, 3).Range.Text = \"36+42=\"\n$t.Cell(17, 4).Range.Text = \"59+28=\"\n$t.Cell(17, 5).Range.Text = \"56+24=\"\n$t.Cell(18, 1).Range.Text = \"74-54=\"\n$t.Cell(18, 2).Range.Text = \"12+41=\"\n$t.Cell(18, 3).Range.Text = \"10+19=\"\n$t.Cell(18, 4).Range.Text = \"85-13=\"\n$t.Cell(18, 5).Range.Text = \"28-3=\"\n$t.Cell(19, 1).Range.Text = \"14+53=\"\n$t.Cell(19, 2).Range.Text = \"80-49=\"\n$t.Cell(19, 3).Range.Text = \"53+0=\"\n$t.Cell(19, 4).Range.Text = \"54-29=\"\n$t.Cell(19, 5).Range.Text = \"35-21=\"\n$t.Cell(20, 1).Range.Text = \"32+64=\"\n$t.Cell(20, 2).Range.Text = \"18+7=\"\n$t.Cell(20, 3).Range.Text = \"0+31=\"\n$t.Cell(20, 4).Range.Text = \"61-25=\"\n$t.Cell(20, 5).Range.Text = \"21+48=\"\n"}
